$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 becomes a filled-in log entry ("Wish list feature developed").
$ws.Rows("11:11").RowHeight = 30

$ws.Range("A11").Value = "Wish list feature developed "
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = "28.97.2024"
$ws.Range("D11").Value = "Created a new page for wishlist, add like functionaility. Liked houses is saved to phone storage."

# Description cell for this row wraps text (matches the other filled rows' style).
$ws.Range("D11").WrapText = $true

# Move the active selection from A15 to A13.
$null = $ws.Range("A13").Select()
